# Simplified keynote creation logic
# Renumber existing keynote codes to start at 0000 within each category,
# and append blank "<Empty>" placeholder rows to each category so that
# new keynotes can be inserted without re-numbering everything again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: renumber the GENERAL (D) category codes down by one
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "D0000"
$ws.Range("A6").Value = "D0001"
$ws.Range("A7").Value = "D0002"
$ws.Range("A8").Value = "D0003"

# Add a brand new blank placeholder row to the GENERAL category, copying
# the formatting of the last existing row in that category (row 8).
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "D0004"
$ws.Range("B9").Value = "<Empty>"
$ws.Range("C9").Value = 0

# ---------------------------------------------------------------------
# Step 2: renumber the EXISTING (E) category codes down by one.
# After the insert above, the E rows moved down one row (was 11/12, now 12/13)
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "E0000"
$ws.Range("A13").Value = "E0001"

# Add three new blank placeholder rows to the EXISTING category, copying
# the formatting of the last existing row in that category (row 13).
$ws.Rows.Item(13).Copy()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "E0002"
$ws.Range("B14").Value = "<Empty>"
$ws.Range("C14").Value = "disabled"

$ws.Range("A15").Value = "E0003"
$ws.Range("B15").Value = "<Empty>"
$ws.Range("C15").Value = "disabled"

$ws.Range("A16").Value = "E0004"
$ws.Range("B16").Value = "<Empty>"
$ws.Range("C16").Value = "disabled"

# ---------------------------------------------------------------------
# Step 3: renumber the NEW (N) category codes down by one.
# These rows have shifted down from 15/16 to 19/20 by now.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "N0000"
$ws.Range("A20").Value = "N0001"

# Add one new blank placeholder row to the NEW category, copying the
# formatting of the last existing row in that category (row 20).
$ws.Rows.Item(20).Copy()
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "N0002"
$ws.Range("B21").Value = "<Empty>"
$ws.Range("C21").Value = "disabled"

$excel.CutCopyMode = $false
